$wb = $excel.ActiveWorkbook

# The workbook currently has two tabs in the order: "2021-Q2", "总计".
# Re-sort the sheet tabs so that "总计" (the summary/totals sheet) comes
# first, followed by "2021-Q2" -- i.e. move "总计" in front of the first
# worksheet. No cell data changes; this only reorders the sheet tabs.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Move($wb.Worksheets.Item(1))

# "2021-Q2" remains the active/selected sheet after the reorder.
$wb.Worksheets.Item("2021-Q2").Activate()
